$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate the newly checked cells (Code Coverage row / Junit row)
# for Bank (D) and Customer (E) columns
$ws.Range("D6").Value = "Y"
$ws.Range("E6").Value = "Y"
$ws.Range("D7").Value = "Y"
$ws.Range("E7").Value = "Y"

# Update the active selection to match the author's final cursor position
$ws.Range("E6").Select()
